# Word COM-interop script
#
# Implements the commit "### 更新 20191222 2033": a new reference link is
# added right after the existing "jianshu.com" hyperlink paragraph. The
# paragraph that used to be an empty spacer (paragraph #7) becomes a
# left-aligned, first-line-indented paragraph that contains a hyperlink to
#   https://www.cnblogs.com/CarpenterLee/p/9558026.html
# formatted the same way as the hyperlink paragraph right above it
# (widowControl on, 420-twip/21-pt first line indent, left justification,
# "微软雅黑" font).

$d = $word.ActiveDocument

# The empty paragraph right after the "jianshu.com" hyperlink paragraph.
$targetIndex = 7
$p = $d.Paragraphs.Item($targetIndex)

# --- paragraph formatting: widowControl / firstLineIndent / left align ---
$p.Format.WidowControl = $true
$p.Format.FirstLineIndent = 21    # points -> 21pt * 20 = 420 twips
$p.Format.Alignment = 0           # wdAlignParagraphLeft

# --- insert the hyperlink text at the (currently empty) paragraph ---
$r = $p.Range
$r.Collapse(0)  # wdCollapseEnd -- land right before the paragraph mark
$url = "https://www.cnblogs.com/CarpenterLee/p/9558026.html"
$d.Hyperlinks.Add($r, $url, [Type]::Missing, [Type]::Missing, $url) | Out-Null

# --- give the new hyperlink run the same "微软雅黑" font used elsewhere ---
$p2 = $d.Paragraphs.Item($targetIndex)
$fr = $p2.Range
$fr.Find.ClearFormatting()
$fr.Find.Replacement.ClearFormatting()
$fr.Find.Replacement.Style = "Default Paragraph Font"
$fr.Find.Replacement.Font.NameAscii = "微软雅黑"
$fr.Find.Replacement.Font.NameFarEast = "微软雅黑"
$fr.Find.Replacement.Font.NameOther = "微软雅黑"
$fr.Find.Execute($url, $true, $false, $false, $false, $false, $true, 1, $false, $url, 2) | Out-Null

Write-Output "Inserted CarpenterLee hyperlink paragraph."
